$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-07 (row 20)
$ws.Range("B20").Value = 6191
$ws.Range("C20").Value = 982
$ws.Range("D20").Value = 5579612
$ws.Range("E20").Value = 901.245679211759
$ws.Range("F20").Value = 6.944204525824849
$ws.Range("G20").Value = 4.24628450106157
$ws.Range("H20").Value = 26.21634760137201
